$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New / changed source-file identifiers used across all three sheets.
# ---------------------------------------------------------------------------
$oldMdFile   = "a6d5f17a-6127-422a-a9c7-e2c1c8202ce2.md"
$newMdFile   = "ba436930-3f19-472c-b819-7d06ea4c6624.md"
$newMdFile2  = "ffff95da5b81-2573-409b-a5f7-2b00479d291b.md"
$cfgFile     = ".localization-config"

$zhXlf       = "ba436930-3f19-472c-b819-7d06ea4c6624.614134b57d741bace7b01fa8ebab4b2f3c7f6b55.zh-cn.xlf"
$deXlf       = "ba436930-3f19-472c-b819-7d06ea4c6624.614134b57d741bace7b01fa8ebab4b2f3c7f6b55.de-de.xlf"

$readyStatus   = "Ready for handoff"
$ignoreStatus  = "Not to be localized"
$epoch         = "0001-01-01 00:00:00"
$zhHandoffTime = "2016-01-20 08:13:18"
$deHandoffTime = "2016-01-20 08:13:29"
$includeWord   = "Include"
$ignoredWord   = "Ignored"

$baseUrl = "https://github.com/OpenLocalizationTest/oltest/blob/8618ad2fdf283a6d8e0cd2a7216d539d46093725/e2e"
$cfgUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/f4803e98e1d7e44f1bd4a1049694b6ab35d00ac1/.localization-config"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview" -- columns A (File Name), B (zh-cn), C (de-de)
#          Adds a new row for the second handed-off file, pushing the
#          ".localization-config" row from row 3 down to row 4.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Delete()

$wsOverview.Range("A2").Value = $newMdFile
$wsOverview.Range("B2").Value = $readyStatus
$wsOverview.Range("C2").Value = $readyStatus

$wsOverview.Range("A3").Value = $newMdFile2
$wsOverview.Range("B3").Value = $readyStatus
$wsOverview.Range("C3").Value = $readyStatus

$wsOverview.Range("A4").Value = $cfgFile
$wsOverview.Range("B4").Value = $ignoreStatus
$wsOverview.Range("C4").Value = $ignoreStatus
$wsOverview.Range("A4").Style = "HyperLink"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "$baseUrl/$newMdFile", "", "", $newMdFile)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "$baseUrl/$newMdFile2", "", "", $newMdFile2)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $cfgUrl, "", "", $cfgFile)

# ---------------------------------------------------------------------------
# Shared helper data for the two per-language detail sheets ("zh-cn", "de-de")
# Columns: A Source File Name, B Status, C Latest Handoff File,
#          D Latest Handoff Datetime, E Latest Target File,
#          F Latest Handback File, G Latest Handback DateTime,
#          H Handoff Reason, I Dependency From
# ---------------------------------------------------------------------------
function Update-LangSheet($ws, $xlfName, $handoffTime) {

    $ws.Hyperlinks.Delete()

    # Row 2: first handed-off source file, now "Ready for handoff" with a
    # real handoff target file + timestamp instead of the old failure state.
    $ws.Range("A2").Value = $newMdFile
    $ws.Range("B2").Value = $readyStatus
    $ws.Range("C2").Value = $xlfName
    $ws.Range("C2").Style = "HyperLink"
    $ws.Range("D2").Value = $handoffTime
    $ws.Range("G2").Value = $epoch
    $ws.Range("H2").Value = $includeWord

    # Row 3 (new): second handed-off source file -- same shape as row 2.
    $ws.Range("A3").Value = $newMdFile2
    $ws.Range("B3").Value = $readyStatus
    $ws.Range("C3").Value = $xlfName
    $ws.Range("C3").Style = "HyperLink"
    $ws.Range("D3").Value = $handoffTime
    $ws.Range("G3").Value = $epoch
    $ws.Range("H3").Value = $includeWord

    # Row 4 (was row 3): the never-localized config file, pushed down.
    $ws.Range("A4").Value = $cfgFile
    $ws.Range("A4").Style = "HyperLink"
    $ws.Range("B4").Value = $ignoreStatus
    $ws.Range("D4").Value = $epoch
    $ws.Range("G4").Value = $epoch
    $ws.Range("H4").Value = $ignoredWord

    $ws.Hyperlinks.Add($ws.Range("A2"), "$baseUrl/$newMdFile", "", "", $newMdFile)
    $ws.Hyperlinks.Add($ws.Range("C2"), "https://example.com/$xlfName", "", "", $xlfName)
    $ws.Hyperlinks.Add($ws.Range("A3"), "$baseUrl/$newMdFile2", "", "", $newMdFile2)
    $ws.Hyperlinks.Add($ws.Range("C3"), "https://example.com/$xlfName", "", "", $xlfName)
    $ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, "", "", $cfgFile)
}

$wsZh = $wb.Worksheets.Item("zh-cn")
Update-LangSheet $wsZh $zhXlf $zhHandoffTime

$wsDe = $wb.Worksheets.Item("de-de")
Update-LangSheet $wsDe $deXlf $deHandoffTime

Write-Output "Report regenerated for handoff"
